$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.087.31"
$ws.Range("E2").Value = "  -1.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.898.80"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.15"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5027"
$ws.Range("E7").Value = "  -0.57%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3901"
$ws.Range("E8").Value = "  -1.38%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09239"
$ws.Range("E9").Value = "  -6.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.128"
$ws.Range("E10").Value = "  -2.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.79"
$ws.Range("E11").Value = "  +0.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.401"
$ws.Range("E12").Value = "  -2.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.81"
$ws.Range("E13").Value = "  -1.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.884.47"
$ws.Range("E14").Value = "  -1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.307"
$ws.Range("E15").Value = "  -3.66%  "

$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.47"
$ws.Range("E17").Value = "  -1.54%  "

$ws.Range("E18").Value = "  -2.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06636"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.89"
$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.231"
$ws.Range("E22").Value = "  -0.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.143.97"
$ws.Range("E23").Value = "  -1.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.47"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.314"
$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.109.42"
$ws.Range("E26").Value = "  -1.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.558"
$ws.Range("E27").Value = "  -6.71%  "

$ws.Range("E28").Value = "  -2.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.17"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.59"
$ws.Range("E30").Value = "  -1.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.083"
$ws.Range("E31").Value = "  -2.28%  "

$ws.Range("E32").Value = "  -0.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.619"
$ws.Range("E33").Value = "  -1.63%  "

$ws.Range("E34").Value = "  -0.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.587"
$ws.Range("E35").Value = "  -3.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06600"
$ws.Range("E36").Value = "  -3.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02399"
$ws.Range("E37").Value = "  -1.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2206"
$ws.Range("E38").Value = "  -1.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.301"
$ws.Range("E39").Value = "  +9.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.226"
$ws.Range("E40").Value = "  -4.19%  "

$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("E42").Value = "  -2.49%  "

$ws.Range("E43").Value = "  -2.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6112"
$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.34"
$ws.Range("E46").Value = "  -2.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.312"
$ws.Range("E47").Value = "  +2.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.687"
$ws.Range("E48").Value = "  +0.64%  "

$ws.Range("E49").Value = "  -2.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.28"
$ws.Range("E50").Value = "  -2.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.196"
$ws.Range("E51").Value = "  -1.28%  "
